$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Add a 4th column ("Stone reference") to the table.
$t.Columns.Add() | Out-Null

# Re-assert every column's width explicitly so Word recalculates the
# tblGrid alongside each cell's tcW (matches target: 3216, 3856, 2896, 2896
# twips == 160.8, 192.8, 144.8, 144.8 points).
$t.Columns.Item(1).Width = 160.8
$t.Columns.Item(2).Width = 192.8
$t.Columns.Item(3).Width = 144.8
$t.Columns.Item(4).Width = 144.8

# Formatted-run XML fragment shared by the three new "Stone ..." cells
# (Times New Roman, sz 24/24, kern 0, ligatures none) identical to the
# formatting already used in the "Created ..." column.
function New-StoneCellXml([string]$text) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:kern w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/><w14:ligatures w14:val="none"/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

$t.Cell(1, 4).Range.InsertXML((New-StoneCellXml "Stone 21")) | Out-Null
$t.Cell(2, 4).Range.InsertXML((New-StoneCellXml "Stone 36, 37, 39")) | Out-Null
$t.Cell(3, 4).Range.InsertXML((New-StoneCellXml "Stone 38, 41")) | Out-Null

# Switch the section to landscape (12240x15840 -> 15840x12240).
$d.PageSetup.Orientation = 1
